$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the staking-strategy input formulas on row 2.
$ws.Range("B2").Formula = "=C2*1.8"
$ws.Range("C2").Formula = "=A2*0.03"

# Add a new (empty) cell at M9 formatted as currency, matching the
# existing "R$" accounting number format used elsewhere in the sheet
# (numFmtId 44) so Excel reuses that format instead of registering a
# new custom one.
$ws.Range("M9").NumberFormat = '_-"R$"\ * #,##0.00_-;\-"R$"\ * #,##0.00_-;_-"R$"\ * "-"??_-;_-@_-'

# Move the current selection to B16, as in the saved workbook.
$ws.Range("B16").Select()
